$wb = $excel.ActiveWorkbook

# --- OFF sheet (Target Depth Data - offense) ---
$wsOff = $wb.Worksheets.Item("OFF")

# Row 2 (H)
$wsOff.Range("B2").Value = 459
$wsOff.Range("C2").Value = 338
$wsOff.Range("D2").Value = 125
$wsOff.Range("E2").Value = 66

# Row 3 (R)
$wsOff.Range("B3").Value = 472
$wsOff.Range("C3").Value = 332
$wsOff.Range("D3").Value = 113
$wsOff.Range("E3").Value = 55
$wsOff.Range("F3").Value = 8

# --- DEF sheet (Target Depth Data - defense) ---
$wsDef = $wb.Worksheets.Item("DEF")

# Row 2 (H)
$wsDef.Range("B2").Value = 510
$wsDef.Range("C2").Value = 351
$wsDef.Range("D2").Value = 130
$wsDef.Range("E2").Value = 48
$wsDef.Range("G2").Value = 5

# Row 3 (R)
$wsDef.Range("B3").Value = 486
$wsDef.Range("C3").Value = 351
$wsDef.Range("D3").Value = 96
$wsDef.Range("E3").Value = 47
